$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the style of the existing last header cell (AC1) onto the new header
# cells so they match the bold/centered/bordered look of the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# New header row (row 1) - columns AD, AE, AF: team win/loss/tie record
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record data for every data row (2 through 39)
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 30).Value = 105  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 57   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
